$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "notes" for the mp3-playback issue (row 6, column D):
# previously a note about browser compatibility, now about the mp3->ima4 file format conversion.
$ws.Range("D6").Value = "FILE FORMAT ISSUE. Must figure out how to convert mp3 into ima4 to save CPU space on iOS"

# The "last and secondLast classes not playing notes" issue (row 11, column C) is now Resolved
# (columns on the ends can now play).
$ws.Range("C11").Value = "Resolved"

# Move the active selection from C19 to C4.
$ws.Range("C4").Select()
